$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 985
$ws.Range("E2").Value = 46200501679
$ws.Range("X2").Value = "DN4127450128190"

# Row 3
$ws.Range("A3").Value = 986
$ws.Range("E3").Value = 46200501680
$ws.Range("X3").Value = "DN4127450128191"

# Row 4
$ws.Range("A4").Value = 987
$ws.Range("E4").Value = 46200501681
$ws.Range("X4").Value = "DN4127450128192"

# Row 5
$ws.Range("A5").Value = 988
$ws.Range("E5").Value = 46200501682
$ws.Range("X5").Value = "DN4127450128193"

# Row 6
$ws.Range("A6").Value = 989
$ws.Range("E6").Value = 46200501683
$ws.Range("X6").Value = "DN4127450128194"
